# ml-la pract. exam 78%
# Adds a new Hungarian/English vocab block (corpus/gyujtemeny, fraction/toredek,
# impute/"betud, berak, "/"impute values for missing data") to Sheet1, right
# after the existing "sparsity/szoras" row, and moves the sheet view/selection
# down to the newly-added rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A29").Value = "corpus"
$ws1.Range("B29").Value = "gyujtemeny"

$ws1.Range("A30").Value = "fraction"
$ws1.Range("B30").Value = "toredek"

$ws1.Range("A31").Value = "impute"
$ws1.Range("B31").Value = "betud, berak, "
$ws1.Range("C31").Value = "impute values for missing data"

# Move the view/selection to match where the author ended up editing.
$ws1.Application.ActiveWindow.ScrollRow = 25
$ws1.Application.ActiveWindow.ScrollColumn = 3
$ws1.Range("C31").Select()
